## Generate Report for Handback
## Updates the zh-cn and de-de localization-status sheets for the
## "4fa3da70-cabe-415e-9c54-00057329f2a7.md" row (row 7): a handback was
## received, but it is not based on the latest source version, so we
## populate the Latest Target File / Latest Handback File / Latest Handback
## DateTime / Error Detail columns accordingly.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fff950fb56cad7c7efa4050fd07eacd4b81e3ad1/e2e/4fa3da70-cabe-415e-9c54-00057329f2a7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af234d41f249b65bea1718b8a8174ba71c88cf40/e2e/4fa3da70-cabe-415e-9c54-00057329f2a7.md."

$targetDisplay = "4fa3da70-cabe-415e-9c54-00057329f2a7.md"

# ---------------------------------------------------------------
# zh-cn sheet (row 7)
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(7, 10).Value = "4fa3da70-cabe-415e-9c54-00057329f2a7.649d27ab56ee89951dfcdd53bb1c7cac58a6c6bf.zh-cn.xlf"
$wsZh.Cells.Item(7, 11).Value = "2016-08-27 22:54:58"
$wsZh.Cells.Item(7, 16).Value = $errorDetail

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item(7, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/649d27ab56ee89951dfcdd53bb1c7cac58a6c6bf/e2e/4fa3da70-cabe-415e-9c54-00057329f2a7.md",
    "",
    "",
    $targetDisplay
) | Out-Null
$wsZh.Cells.Item(7, 9).Style = "HyperLink"

# ---------------------------------------------------------------
# de-de sheet (row 7)
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(7, 10).Value = "4fa3da70-cabe-415e-9c54-00057329f2a7.649d27ab56ee89951dfcdd53bb1c7cac58a6c6bf.de-de.xlf"
$wsDe.Cells.Item(7, 11).Value = "2016-08-27 22:55:11"
$wsDe.Cells.Item(7, 16).Value = $errorDetail

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item(7, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/649d27ab56ee89951dfcdd53bb1c7cac58a6c6bf/e2e/4fa3da70-cabe-415e-9c54-00057329f2a7.md",
    "",
    "",
    $targetDisplay
) | Out-Null
$wsDe.Cells.Item(7, 9).Style = "HyperLink"

Write-Host "Report generated for handback row."
